# The deck ships two DrawingML theme parts: ppt/theme/theme1.xml (stock
# "Office Theme" colours) and ppt/theme/theme2.xml (the "Integral" theme
# that the single slide master / the presentation actually use). The
# target edit swaps their bodies so the live design's colour scheme
# reverts to the default Office palette (and the dormant part picks up
# the Integral palette). Font scheme / effect scheme are identical
# between the two themes, so only the 12 theme colours actually change
# for the part that PowerPoint's object model can reach (the slide
# master's live Theme == ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Office Theme colour scheme, in clrScheme document order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), expressed as
# OLE RGB() integers (0x00BBGGRR).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}

# Font scheme is already Arial/Arial on both themes, but set it
# explicitly so the live theme's font scheme matches the target too.
$fontScheme = $theme.ThemeElements.FontScheme
$fontScheme.MajorFont.Latin = "Arial"
$fontScheme.MinorFont.Latin = "Arial"
